$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item("testSheetNames"))
$ws.Name = "testSheetYears"

$ws.Range("A1").Value = "TABID"
$ws.Range("B1").Value = "NOYEARS"
$ws.Range("A2").Value = "P509633.conll"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "mu {d}nin-urta u3 ma-na-ba-al-te-el in-pa3 igi ia-ah-zi-be2-el3 dumu isz-me-{d}suen igi szu-la-ra-ak igi na-nu-szu-um dumu za-x-la-nu igi ku-nu-um dumu nu-ur2-ia "
$ws.Range("D2").Value = "mu e2 {d}nin-urta ma-na-ba-al-te-el mu-na-du3 dumu isz-me-{d}suen igi szu-la-ra-ak dumu x-x-x-x igi ku-nu-um dumu nu-ur2-ia ... x ... "
$ws.Range("A3").Value = "P509634.conll"
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = "P509636.conll"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "mu {gesz}gu-za za3-bi-us2 3(disz)-a-bi {d}en-lil2 {d}nin-urta {d}nin-nibru{ki} ma-na-ba-al-te-el {disz}{d}suen-a-hu-szu u3 sza-at-{d}gibil6 ha-la e2 ad-da-ni nig2 na-me ugu za-za-kum li-bi2-in-tuku igi a-mur-i-lu-su2 igi {d}suen-en-nam szitim igi bur-ia dumu sa-li igi {d}nanna-zi-sza3-gal2 iti gan-gan-e3 "
$ws.Range("D4").Value = "mu {gesz}gu-za za3-bi-us2 3(disz)-a-bi {d}en-lil2 {d}nin-urta {d}nin-nibru{ki} ma-na-dim2 {d}suen-a-bu-szu dumu za-pa-ti "
$ws.Range("A5").Value = "P509663.conll"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "mu ma-ma-ba-al-te-el 6(disz) gin2 ku3-babbar sa10 2(asz) gur sze szu ti-a a-na-gum "
$ws.Range("D5").Value = "mu mu-na-ba-al-te-el 1(u) 3(disz) 1/2(disz) gin2 1(disz)-kam 1(u) 3(disz) 1/3(disz) gin2 2(disz)-kam 1/3(disz) ma-na 8(disz) igi-6(disz)-gal2 szu-la2 "
$ws.Range("E5").Value = "mu ma-na-ba-al-te-el "
$ws.Range("A6").Value = "P509664.conll"
$ws.Range("B6").Value = 0
$ws.Range("A7").Value = "P509668.conll"
$ws.Range("B7").Value = 0
$ws.Range("A8").Value = "P509669.conll"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "mu ... 4(disz) gu4-ab2 6(disz) ab2 mah2 2(u) 7(disz) gu4 amar-ga 1(disz) amar nig2 e2 nam-ra-tum nig2-szu la-lum giri3 x-x-x-x i3-gen-ne-en iti sze-sag11-ku5 u4 4(disz)-kam "
$ws.Range("D8").Value = "mu GISZ ... la-lu-um dumu sa-al-lum "
$ws.Range("A9").Value = "P509670.conll"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "mu ... 4(disz) ab2 mah2 6(disz) gu4-ab2 2(u) 2(disz) ab2 amar-ga 6(disz) gu4 amar-ga 3(u) 8(disz) gu4 ab2 amar-hi-a nig2-szu {d}utu-sipa "
$ws.Range("A10").Value = "P509671.conll"
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = "mu 1(disz) 2(disz) gu4 "
$ws.Range("D10").Value = "mu 2(disz) 1(disz) gu4 "
$ws.Range("E10").Value = "mu 1(disz) 1(disz) gu4 amar-ga 1(disz) ab2 amar-ga 1(u) 1(disz) gu4 ab2 amar dingir-na-sir 2(disz) ab2 mah2 2(disz) ab2 "
$ws.Range("F10").Value = "mu 1(disz) 1(disz) ab2 amar-ga 1(disz) gu4 amar-ga 6(disz) er3-ra-ha-bi-it 4(disz) ab2 mah2 2(disz) ab2 amar-ga 1(disz) gu4 amar-ga n ab2 amar-ga n ab2 "
$ws.Range("G10").Value = "mu n n gu4 "
$ws.Range("H10").Value = "mu 2(disz) 1(disz) gu4 "
$ws.Range("I10").Value = "mu 1(disz) 1(disz) gu4 amar-ga 1(disz) ab2 amar-ga 8(disz) a-hu-um-ma 1(disz) ab2 mah2 1(disz) gu4 "
$ws.Range("J10").Value = "mu 3(disz) 1(disz) gu4 "
$ws.Range("K10").Value = "mu 1(disz) 1(disz) ab2 amar-ga 4(disz) {d}suen-na-si-ir 2(disz) ab2 mah2 2(disz) ab2 "
$ws.Range("L10").Value = "mu 2(disz) 1(disz) gu4 "
$ws.Range("M10").Value = "mu n 5(disz) sza3 2(u) a-ra2 1(disz)-kam 5(disz) ab2 mah2 2(disz) ab2 "
$ws.Range("N10").Value = "mu 2(disz) 2(disz) ab2 "
$ws.Range("O10").Value = "mu 1(disz) 2(disz) gu4 "
$ws.Range("P10").Value = "mu 1(disz) 3(disz) gu4 "
$ws.Range("Q10").Value = "mu 2(disz) 3(disz) gu4 "
$ws.Range("R10").Value = "mu 3(disz) 1(u) 7(disz) x ... nig2-szu {d}suen-x-x "

$ws.Range("L10").Select() | Out-Null
